$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "caseDetailQuery" column (column C) entirely, which contained the
# query with the hardcoded/dynamic case id. This shifts the dbExcel / WebExcel
# columns (old D, E) left into C, D.
$ws.Range("C1:C2").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftToLeft) | Out-Null

# Update the active selection to reflect the post-edit state.
$ws.Range("J2").Select() | Out-Null
